$wb = $excel.ActiveWorkbook

# Insert two new sheets ("Revenue statistics", "Currency rate") right after
# "MAU statistics" and right before "Step-by-step statistics".
$mauSheet = $wb.Worksheets.Item("MAU statistics")

$revenue = $wb.Worksheets.Add($null, $mauSheet)
$revenue.Name = "Revenue statistics"
$revenue.Range("A1").Value = "Day"
$revenue.Range("B1").Value = "Revenue, `$"
$revenue.Cells.Item(2, 1).Formula = "=""01.01.2018"""
$revenue.Cells.Item(2, 1).Copy()
$revenue.Cells.Item(2, 1).PasteSpecial(-4163)
$revenue.Range("B2").Value = 11054

$currency = $wb.Worksheets.Add($null, $revenue)
$currency.Name = "Currency rate"
$currency.Range("A1").Value = "Day"
$currency.Range("B1").Value = "Rate, `$ / curr"
$currency.Cells.Item(2, 1).Formula = "=""01.01.2018"""
$currency.Cells.Item(2, 1).Copy()
$currency.Cells.Item(2, 1).PasteSpecial(-4163)
$currency.Range("B2").Value = 0.08668104293275829
